$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 with the latest quotation date and values
$ws.Range("A4").Value = 45908
$ws.Range("A4").Style = $ws.Range("A2").Style
$ws.Range("A4").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B4").Value = "20,7525"
$ws.Range("C4").Value = "14,6423"
$ws.Range("D4").Value = "14,7257"
$ws.Range("E4").Value = "14,7257"
